# Daily attendance processing - 2025-12-28 19:28:52
#
# Normalizes the "Recorded By" list in column G of the session-analysis
# sheet: whenever a cell holds a comma-separated list of recorders whose
# last entry is not "System", that last entry is moved to the front of
# the list (the rest keep their relative order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column A (data starts on row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $text = [string]$cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text.Split(",")
    for ($i = 0; $i -lt $parts.Count; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Count -gt 1 -and $parts[$parts.Count - 1] -ne "System") {
        $newParts = @($parts[$parts.Count - 1]) + $parts[0..($parts.Count - 2)]
        $newVal = [string]::Join(", ", $newParts)
        $cell.Value = $newVal
    }
}

Write-Host "Recorded By column normalized through row $lastRow"
